$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows: reassign "In-charge" from Admin to Staff ---
# (B9 currently has no quote-prefix style, plain value assignment keeps its style)
$ws.Range("B9").Value = "Staff"

# B22 currently uses the quote-prefixed style variant (s=6). A plain .Value
# assignment collapses that to the non quote-prefixed style (s=5), so restore
# the original formatting afterwards by pasting formats from a sibling cell
# that still carries the same style (B23).
$ws.Range("B22").Value = "Staff"
$ws.Range("B23").Copy() | Out-Null
$ws.Range("B22").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- Append 2 new backlog entries (rows 39-41) ---
# Copy the formatting of the last existing row (38) down into the new rows
# first, so the new cells inherit the same styles (incl. the empty D/E/F/G
# placeholder cells), then fill in the text.
$ws.Range("A38:G38").Copy() | Out-Null
$ws.Range("A39:G41").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("A39").Value = "Move physical document "
$ws.Range("B39").Value = "Staff"
$ws.Range("C39").Value = "Simple"

$ws.Range("A40").Value = "Auhorization physical document"
$ws.Range("B40").Value = "Admin"
$ws.Range("C40").Value = "Complex"

$ws.Range("A41").Value = "Authorization digital document"
$ws.Range("B41").Value = "Admin"
$ws.Range("C41").Value = "Complex"

$excel.CutCopyMode = 0

# --- Extend the data validation (dropdown list) ranges to cover the new rows ---
$ws.Range("D9:D38").Validation.Delete() | Out-Null
$ws.Range("G9:G38").Validation.Delete() | Out-Null
$ws.Range("C9:C38").Validation.Delete() | Out-Null

$ws.Range("D9:D41").Validation.Add(3, 1, 1, '"Iteration 1, Iteration 2, Iteration 3"') | Out-Null
$ws.Range("G9:G41").Validation.Add(3, 1, 1, '"Iteration 1, Iteration 2, Iteration 3, Final"') | Out-Null
$ws.Range("C9:C41").Validation.Add(3, 1, 1, '"Simple, Medium, Complex"') | Out-Null

# --- Update the view: scroll down a couple of rows and move the active
# selection to the first blank row below the new data ---
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("A42").Select() | Out-Null
